$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark that currently sits in the very
#    first paragraph of the document.
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# ------------------------------------------------------------------
# 2) The last paragraph of the body is empty (it only carries the
#    paragraph-mark formatting used for the "End of Document." run).
#    Build the exact run we need - same rFonts / color / sz / szCs /
#    lang as the paragraph mark - and splice it (plus a fresh
#    "_GoBack" bookmark right after it) into that paragraph using
#    InsertXML so the run-level formatting round-trips exactly.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertionRange = $lastPara.Range

$paragraphProps = '<w:pPr>' `
    + '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' `
    + '<w:spacing w:before="196" w:after="196" w:line="346" w:lineRule="atLeast"/>' `
    + '<w:outlineLvl w:val="2"/>' `
    + '<w:rPr>' `
    + '<w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/>' `
    + '<w:color w:val="0070C0"/>' `
    + '<w:sz w:val="29"/>' `
    + '<w:szCs w:val="29"/>' `
    + '<w:lang w:val="en-US" w:eastAsia="fr-FR"/>' `
    + '</w:rPr>' `
    + '</w:pPr>'

$newRun = '<w:r>' `
    + '<w:rPr>' `
    + '<w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/>' `
    + '<w:color w:val="0070C0"/>' `
    + '<w:sz w:val="29"/>' `
    + '<w:szCs w:val="29"/>' `
    + '<w:lang w:val="en-US" w:eastAsia="fr-FR"/>' `
    + '</w:rPr>' `
    + '<w:t>End of Document.</w:t>' `
    + '</w:r>'

$newBookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$newParagraphXml = '<w:p>' + $paragraphProps + $newRun + $newBookmark + '</w:p>'

$openXmlPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData>' `
    + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:body>' + $newParagraphXml + '</w:body>' `
    + '</w:document>' `
    + '</pkg:xmlData></pkg:part></pkg:package>'

# InsertXML splices a whole new paragraph in just before the
# (still-empty) original paragraph mark, so the document temporarily
# gains one extra paragraph.
$insertionRange.InsertXML($openXmlPackage)

# ------------------------------------------------------------------
# 3) Merge the newly inserted paragraph back together with the
#    original (now trailing, still-empty) paragraph mark by deleting
#    the paragraph break between them.
# ------------------------------------------------------------------
$newLastIndex = $d.Paragraphs.Count - 1
$mergedPara = $d.Paragraphs.Item($newLastIndex)
$paraMarkRange = $d.Range($mergedPara.Range.End - 1, $mergedPara.Range.End)
$paraMarkRange.Delete()

Write-Output ("Paragraphs: {0}" -f $d.Paragraphs.Count)
Write-Output ("Final paragraph text: [{0}]" -f $d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)
